# "Big update for 02Jan" - update confirmed-death counts for several existing
# days near the end of the series, and append six new days of data
# (2020-12-26 .. 2020-12-31 / Excel serials 44191..44196).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Corrections to column B (Confirmed Deaths) on already-existing rows.
#    Columns C (running total), E (probable running total) and F (7-day
#    average) are formulas and recalculate automatically.
# ---------------------------------------------------------------------------
$bEdits = @{
    248 = 32
    279 = 48
    284 = 56
    285 = 56
    287 = 58
    289 = 69
    290 = 48
    291 = 79
    292 = 64
}
foreach ($row in $bEdits.Keys) {
    $ws.Cells.Item($row, 2).Value = $bEdits[$row]
}

# ---------------------------------------------------------------------------
# 2) Corrections to column D (Probable Deaths) on already-existing rows.
# ---------------------------------------------------------------------------
$dEdits = @{
    277 = 3
    285 = 1
    292 = 1
}
foreach ($row in $dEdits.Keys) {
    $ws.Cells.Item($row, 4).Value = $dEdits[$row]
}

# ---------------------------------------------------------------------------
# 3) Append six new rows (293-298) for 2020-12-26 .. 2020-12-31.
#    Copy the date format from A292 first so the new date cells get the
#    same m/d/yyyy display as the rest of column A.
# ---------------------------------------------------------------------------
$ws.Range("A292").Copy() | Out-Null
$ws.Range("A293:A298").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row = 293; Date = 44191; B = 58; D = 1 },
    @{ Row = 294; Date = 44192; B = 58; D = 2 },
    @{ Row = 295; Date = 44193; B = 72; D = 0 },
    @{ Row = 296; Date = 44194; B = 70; D = 3 },
    @{ Row = 297; Date = 44195; B = 41; D = 0 },
    @{ Row = 298; Date = 44196; B = 3;  D = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Formula = "=B$row+C$($row-1)"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Formula = "=D$row+E$($row-1)"
    $ws.Cells.Item($row, 6).Formula = "=AVERAGE(B$($row-6):B$row)"
}

# ---------------------------------------------------------------------------
# 4) Leave the final selection on I299, matching where the author's cursor
#    ended up after entering the new data.
# ---------------------------------------------------------------------------
$ws.Range("I299").Select() | Out-Null
